$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "36.471.20"
$ws.Range("E2").Value = "  +1.52%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.946.52"
$ws.Range("E3").Value = "  -0.57%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.12%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'243.55"
$ws.Range("E5").Value = "  +0.80%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -1.48%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.05%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "'0.366"
$ws.Range("E9").Value = "  -0.60%  "

# Row 10 - OKB
$ws.Range("D10").Value = "'55.68"
$ws.Range("E10").Value = "  -0.59%  "

# Row 11 - Dogecoin
$ws.Range("D11").Value = "'0.0834"
$ws.Range("E11").Value = "  +4.54%  "

# Row 12 - TRON
$ws.Range("D12").Value = "'0.104"
$ws.Range("E12").Value = "  +1.02%  "

# Row 13 - was Avalanche, now Polygon
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "'0.823"
$ws.Range("E13").Value = "  -3.82%  "

# Row 14 - was Polygon, now Avalanche
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "'21.56"
$ws.Range("E14").Value = "  -2.07%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.232.77"
$ws.Range("E15").Value = "  -0.46%  "

# Row 16 - Chainlink
$ws.Range("D16").Value = "'13.56"
$ws.Range("E16").Value = "  -3.09%  "

# Row 17 - Polkadot
$ws.Range("E17").Value = "  -2.91%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "1.955.84"
$ws.Range("E18").Value = "  -0.43%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "36.393.85"
$ws.Range("E19").Value = "  +1.64%  "

# Row 20 - Litecoin
$ws.Range("D20").Value = "'69.79"
$ws.Range("E20").Value = "  -1.63%  "

# Row 21 - ShibaInu
$ws.Range("E21").Value = "  +1.32%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "'229.81"
$ws.Range("E22").Value = "  -3.26%  "

# Row 23 - Uniswap
$ws.Range("E23").Value = "  -2.21%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  +0.00%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  -3.49%  "

# Row 26 - Toncoin
$ws.Range("E26").Value = "  -0.04%  "

# Row 27 - Cosmos
$ws.Range("D27").Value = "'9.21"
$ws.Range("E27").Value = "  -6.00%  "

# Row 28 - Monero
$ws.Range("D28").Value = "'162.07"
$ws.Range("E28").Value = "  +1.78%  "

# Row 29 - EthereumClassic
$ws.Range("D29").Value = "'19.42"
$ws.Range("E29").Value = "  -1.83%  "

# Row 30 - Kaspa
$ws.Range("D30").Value = "'0.127"
$ws.Range("E30").Value = "  -1.86%  "

# Row 31 - Stellar
$ws.Range("E31").Value = "  -1.27%  "

# Row 32 - ImmutableX
$ws.Range("E32").Value = "  +1.33%  "

# Row 33 - Filecoin
$ws.Range("E33").Value = "  -3.26%  "

# Row 34 - Hedera
$ws.Range("D34").Value = "'0.0627"
$ws.Range("E34").Value = "  +1.51%  "

# Row 35 - InternetComputer(DFINITY)
$ws.Range("D35").Value = "'4.29"
$ws.Range("E35").Value = "  -2.31%  "

# Row 36 - THORChain
$ws.Range("D36").Value = "'6.21"
$ws.Range("E36").Value = "  -0.86%  "

# Row 37 - BinanceUSD
$ws.Range("E37").Value = "  +0.02%  "

# Row 38 - WEMIXToken
$ws.Range("E38").Value = "  -3.23%  "

# Row 39 - LidoDAOToken
$ws.Range("E39").Value = "  -6.02%  "

# Row 40 - RenderToken
$ws.Range("D40").Value = "'3.03"
$ws.Range("E40").Value = "  -2.03%  "

# Row 41 - Cronos
$ws.Range("D41").Value = "'0.0986"
$ws.Range("E41").Value = "  +0.29%  "

# Row 42 - HuobiToken
$ws.Range("D42").Value = "'2.87"
$ws.Range("E42").Value = "  +2.43%  "

# Row 43 - TrustWalletToken
$ws.Range("D43").Value = "'1.18"
$ws.Range("E43").Value = "  -3.49%  "

# Row 44 - VeChain
$ws.Range("E44").Value = "  -0.83%  "

# Row 45 - InjectiveProtocol
$ws.Range("D45").Value = "'16.10"
$ws.Range("E45").Value = "  +0.01%  "

# Row 46 - Maker
$ws.Range("D46").Value = "1.350.10"
$ws.Range("E46").Value = "  +1.11%  "

# Row 47 - ARBITRUM
$ws.Range("E47").Value = "  -4.65%  "

# Row 48 - Aave
$ws.Range("D48").Value = "'87.76"
$ws.Range("E48").Value = "  -4.57%  "

# Row 49 - FraxShare
$ws.Range("D49").Value = "'7.17"
$ws.Range("E49").Value = "  -4.28%  "

# Row 50 - MXToken
$ws.Range("E50").Value = "  +1.82%  "

# Row 51 - MultiversX
$ws.Range("E51").Value = "  +4.04%  "
